$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These columns hold numeric-looking values that are stored as TEXT in the
# source workbook (t="str"). Force the Text number format before assigning
# so Excel keeps them as text instead of auto-converting to numbers.
$cells = @("C2", "D2", "F2", "C3", "D3", "F3")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: runs 28 -> 51, balls 18 -> 30, sixes 1 -> 2 (fours stays 4)
$ws.Range("C2").Value = "51"
$ws.Range("D2").Value = "30"
$ws.Range("F2").Value = "2"

# Row 3: runs 51 -> 28, balls 30 -> 18, sixes 2 -> 1 (fours stays 4)
$ws.Range("C3").Value = "28"
$ws.Range("D3").Value = "18"
$ws.Range("F3").Value = "1"
